$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block inserted between the existing "bikeModel" block (rows 6-21) and
# the next block (row 25 onward): an extra Item_Attribute_Name/content_location
# pair at rows 22-23.
$ws.Range("A22").Value = "Item_Attribute_Name"
$ws.Range("B22").Value = "bikeModel"
$ws.Range("B22").HorizontalAlignment = -4131

$ws.Range("A23").Value = "content_location"

# New block appended after the last existing block (row 38): another
# Item_Attribute_Name/content_location pair at rows 40-41 (row 39 left blank,
# mirroring the blank-row separators used between the other blocks).
$ws.Range("A40").Value = "Item_Attribute_Name"
$ws.Range("B40").Value = "price"

$ws.Range("A41").Value = "content_location"

# Restore the view state: scrolled so row 15 is the top-left visible row,
# with A23 as the active selected cell (matches the edited sheetView).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A23").Select() | Out-Null
